# Daily attendance processing - 2025-12-03 14:56:59
# Normalizes the "Recorded By" (column G) value ordering on the
# "Session Analysis Results" sheet. The set of recorder names stays the
# same for every row; only the order in which they are listed changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of old "Recorded By" text -> new "Recorded By" text.
# Built from the actual distinct values found in column G; any value not
# present in this map is left untouched.
$valueMap = @{
    "backup@backdoor.com, system, System" = "backup@backdoor.com, System, system"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $current = $cell.Value()
    if ($null -ne $current -and $valueMap.ContainsKey($current)) {
        $cell.Value = $valueMap[$current]
    }
}
